$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Medulloblastoma")

$ws.Range("B2").Value = 0.099269066823374
$ws.Range("B3").Value = 0.647110170368091
$ws.Range("B4").Value = 0.976232230683878
$ws.Range("B5").Value = 0.950212628821588
$ws.Range("B6").Value = 0.410362840774718
